$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 148.1
$ws.Range("I18").Value = 148.1
$ws.Range("K18").Value = 148.1
$ws.Range("M18").Value = 135.9
$ws.Range("H69").Value = 3600.1738
$ws.Range("I69").Value = 3019.1667
$ws.Range("J69").Value = 4234
$ws.Range("K69").Value = 9057.500100000001
$ws.Range("L69").Value = 12702
$ws.Range("M69").Value = -8183.500100000001
$ws.Range("N69").Value = -14450
$ws.Range("H72").Value = 3600.1738
$ws.Range("I72").Value = 3019.1667
$ws.Range("J72").Value = 4234
$ws.Range("K72").Value = 27172.5003
$ws.Range("L72").Value = 38106
$ws.Range("M72").Value = -22804.5003
$ws.Range("N72").Value = -46842
$ws.Range("I74").Value = 1575
$ws.Range("J74").Value = 3902.2222
$ws.Range("K74").Value = 1575
$ws.Range("L74").Value = 3902.2222
$ws.Range("M74").Value = -639
$ws.Range("N74").Value = -5774.2222
$ws.Range("I77").Value = 1575
$ws.Range("J77").Value = 3902.2222
$ws.Range("K77").Value = 7875
$ws.Range("L77").Value = 19511.111
$ws.Range("M77").Value = -3195
$ws.Range("N77").Value = -28871.111
$ws.Range("H87").Value = 12805.857
$ws.Range("J87").Value = 12805.857
$ws.Range("L87").Value = 12805.857
$ws.Range("N87").Value = -15301.857
$ws.Range("H90").Value = 12805.857
$ws.Range("J90").Value = 12805.857
$ws.Range("L90").Value = 38417.571
$ws.Range("N90").Value = -50897.571
$ws.Range("H112").Value = 1107.0426
$ws.Range("J112").Value = 1118.4667
$ws.Range("L112").Value = 3355.4001
$ws.Range("N112").Value = -5571.4001
$ws.Range("H113").Value = 4883.0205
$ws.Range("I113").Value = 2327.8386
$ws.Range("J113").Value = 9283.611000000001
$ws.Range("K113").Value = 2327.8386
$ws.Range("L113").Value = 9283.611000000001
$ws.Range("M113").Value = 926.1614
$ws.Range("N113").Value = -15791.611
$ws.Range("H137").Value = 1884.5
$ws.Range("I137").Value = 1884.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5653.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3103.5
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5192.43
$ws.Range("I32").Value = 4285.1816
$ws.Range("K32").Value = 4285.1816
$ws.Range("M32").Value = -3998.1816
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H45").Value = 1674.4117
$ws.Range("I45").Value = 1604.6428
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1604.6428
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1227.6428
$ws.Range("N45").Value = -2754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H127").Value = 40581
$ws.Range("J127").Value = 40581
$ws.Range("L127").Value = 40581
$ws.Range("N127").Value = -50501
$ws.Range("H134").Value = 1052.2344
$ws.Range("I134").Value = 907.7895
$ws.Range("K134").Value = 2723.3685
$ws.Range("M134").Value = -188.3685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1754.9694
$ws.Range("J31").Value = 3632.9333
$ws.Range("L31").Value = 3632.9333
$ws.Range("N31").Value = -4222.933300000001
$ws.Range("H34").Value = 1754.9694
$ws.Range("J34").Value = 3632.9333
$ws.Range("L34").Value = 3632.9333
$ws.Range("N34").Value = -4036.9333
$ws.Range("H98").Value = 32000
$ws.Range("J98").Value = 32000
$ws.Range("L98").Value = 32000
$ws.Range("N98").Value = -36492
$ws.Range("H99").Value = 14712.125
$ws.Range("I99").Value = 2116.1667
$ws.Range("J99").Value = 52500
$ws.Range("K99").Value = 2116.1667
$ws.Range("L99").Value = 52500
$ws.Range("M99").Value = -618.1667000000002
$ws.Range("N99").Value = -55496
$ws.Range("H122").Value = 3832.75
$ws.Range("I122").Value = 3606
$ws.Range("J122").Value = 4059.5
$ws.Range("K122").Value = 10818
$ws.Range("L122").Value = 12178.5
$ws.Range("M122").Value = -8368
$ws.Range("N122").Value = -17078.5
$ws.Range("H126").Value = 14712.125
$ws.Range("I126").Value = 2116.1667
$ws.Range("J126").Value = 52500
$ws.Range("K126").Value = 6348.500100000001
$ws.Range("L126").Value = 157500
$ws.Range("M126").Value = -3878.500100000001
$ws.Range("N126").Value = -162440
$ws.Range("H134").Value = 2102.08
$ws.Range("I134").Value = 2205.182
$ws.Range("J134").Value = 1346
$ws.Range("K134").Value = 6615.545999999999
$ws.Range("L134").Value = 4038
$ws.Range("M134").Value = -4080.545999999999
$ws.Range("N134").Value = -9108

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10.166667
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10.2
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 61.2
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -287.2
$ws.Range("H35").Value = 2299.75
$ws.Range("J35").Value = 2999.6667
$ws.Range("L35").Value = 8999.000100000001
$ws.Range("N35").Value = -9575.000100000001
$ws.Range("H122").Value = 1268.6154
$ws.Range("I122").Value = 485
$ws.Range("J122").Value = 4559.8
$ws.Range("K122").Value = 4365
$ws.Range("L122").Value = 41038.2
$ws.Range("M122").Value = -1915
$ws.Range("N122").Value = -45938.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H75").Value = 22500
$ws.Range("I75").Value = 15000
$ws.Range("K75").Value = 15000
$ws.Range("M75").Value = -14126
$ws.Range("H78").Value = 22500
$ws.Range("I78").Value = 15000
$ws.Range("K78").Value = 45000
$ws.Range("M78").Value = -40632
$ws.Range("H126").Value = 1926.9788
$ws.Range("I126").Value = 1709.4667
$ws.Range("J126").Value = 2310.8235
$ws.Range("K126").Value = 5128.4001
$ws.Range("L126").Value = 6932.470499999999
$ws.Range("M126").Value = -2658.4001
$ws.Range("N126").Value = -11872.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1868.7142
$ws.Range("I22").Value = 1033.6666
$ws.Range("J22").Value = 2495
$ws.Range("K22").Value = 1033.6666
$ws.Range("L22").Value = 2495
$ws.Range("M22").Value = -738.6666
$ws.Range("N22").Value = -3085
$ws.Range("H27").Value = 1868.7142
$ws.Range("I27").Value = 1033.6666
$ws.Range("J27").Value = 2495
$ws.Range("K27").Value = 1033.6666
$ws.Range("L27").Value = 2495
$ws.Range("M27").Value = -926.6666
$ws.Range("N27").Value = -2709

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 23333.334
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 23333.334
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 23333.334
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -23919.334
$ws.Range("H80").Value = 77777
$ws.Range("J80").Value = 77777
$ws.Range("L80").Value = 77777
$ws.Range("N80").Value = -79773
$ws.Range("H83").Value = 77777
$ws.Range("J83").Value = 77777
$ws.Range("L83").Value = 233331
$ws.Range("N83").Value = -243315
$ws.Range("H126").Value = 84550.414
$ws.Range("I126").Value = 100770
$ws.Range("J126").Value = 3452.5
$ws.Range("K126").Value = 302310
$ws.Range("L126").Value = 10357.5
$ws.Range("M126").Value = -299840
$ws.Range("N126").Value = -15297.5
$ws.Range("I136").Value = 500698.16
$ws.Range("J136").Value = 6415.952
$ws.Range("K136").Value = 1502094.48
$ws.Range("L136").Value = 19247.856
$ws.Range("M136").Value = -1499544.48
$ws.Range("N136").Value = -24347.856
